$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.034777666666667
$ws.Range("H2").Value = 9.104333
$ws.Range("I2").Value = 0.2502264227183869
$ws.Range("J2").Value = 0.2502264227183869
$ws.Range("M2").Value = 0.3045636666666667
$ws.Range("N2").Value = 0.913691
$ws.Range("Q2").Value = 0.9242830136781112
$ws.Range("R2").Value = 8.318547123103
$ws.Range("S2").Value = 0.2502264227183869
$ws.Range("T2").Value = 0.2502264227183869

# Row 3
$ws.Range("I3").Value = 0.4835045831069426
$ws.Range("J3").Value = 0.4835045831069426
$ws.Range("M3").Value = 0.3045636666666667
$ws.Range("N3").Value = 0.913691
$ws.Range("Q3").Value = 1.785962762630444
$ws.Range("R3").Value = 16.073664863674
$ws.Range("S3").Value = 0.4835045831069426
$ws.Range("T3").Value = 0.4835045831069426

# Row 4
$ws.Range("G4").Value = 2.564975
$ws.Range("H4").Value = 7.694925
$ws.Range("I4").Value = 0.2114897989601526
$ws.Range("J4").Value = 0.2114897989601526
$ws.Range("M4").Value = 0.3045636666666667
$ws.Range("N4").Value = 0.913691
$ws.Range("Q4").Value = 0.7811981909083333
$ws.Range("R4").Value = 7.030783718175
$ws.Range("S4").Value = 0.2114897989601526
$ws.Range("T4").Value = 0.2114897989601526

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6643690000000001
$ws.Range("H5").Value = 1.993107
$ws.Range("I5").Value = 0.05477919521451775
$ws.Range("J5").Value = 0.05477919521451775
$ws.Range("M5").Value = 0.3045636666666667
$ws.Range("N5").Value = 0.913691
$ws.Range("Q5").Value = 0.2023426586596667
$ws.Range("R5").Value = 1.821083927937
$ws.Range("S5").Value = 0.05477919521451775
$ws.Range("T5").Value = 0.05477919521451775
